# Auto-generated Excel COM-interop script
# Applies updated TPM-derived values to Col2a1-Ddr1 LR-pairs sheet (rows 2-26)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.046494
$ws.Range("H2").Value = 0.139482
$ws.Range("I2").Value = 0.05587126560541624
$ws.Range("J2").Value = 0.05587126560541625
$ws.Range("M2").Value = 0.7592313333333333
$ws.Range("N2").Value = 2.277694
$ws.Range("O2").Value = 0.05311309006881704
$ws.Range("P2").Value = 0.05311309006881704
$ws.Range("Q2").Value = 0.035299701612
$ws.Range("R2").Value = 0.317697314508
$ws.Range("S2").Value = 0.002967495562359272
$ws.Range("T2").Value = 0.002967495562359273

# Row 3
$ws.Range("G3").Value = 0.046494
$ws.Range("H3").Value = 0.139482
$ws.Range("I3").Value = 0.05587126560541624
$ws.Range("J3").Value = 0.05587126560541625
$ws.Range("O3").Value = 0.1721303313829956
$ws.Range("P3").Value = 0.1721303313829956
$ws.Range("Q3").Value = 0.114400222776
$ws.Range("R3").Value = 1.029602004984
$ws.Range("S3").Value = 0.00961713946344766
$ws.Range("T3").Value = 0.009617139463447662

# Row 4
$ws.Range("G4").Value = 0.046494
$ws.Range("H4").Value = 0.139482
$ws.Range("I4").Value = 0.05587126560541624
$ws.Range("J4").Value = 0.05587126560541625
$ws.Range("M4").Value = 1.018760333333333
$ws.Range("N4").Value = 3.056281
$ws.Range("O4").Value = 0.07126880433834142
$ws.Range("P4").Value = 0.07126880433834142
$ws.Range("Q4").Value = 0.047366242938
$ws.Range("R4").Value = 0.426296186442
$ws.Range("S4").Value = 0.003981878296567915
$ws.Range("T4").Value = 0.003981878296567915

# Row 5
$ws.Range("G5").Value = 0.046494
$ws.Range("H5").Value = 0.139482
$ws.Range("I5").Value = 0.05587126560541624
$ws.Range("J5").Value = 0.05587126560541625
$ws.Range("M5").Value = 9.893154666666666
$ws.Range("N5").Value = 29.679464
$ws.Range("O5").Value = 0.6920894749804902
$ws.Range("P5").Value = 0.6920894749804903
$ws.Range("Q5").Value = 0.459972333072
$ws.Range("R5").Value = 4.139750997648
$ws.Range("S5").Value = 0.03866791487934805
$ws.Range("T5").Value = 0.03866791487934806

# Row 6
$ws.Range("G6").Value = 0.046494
$ws.Range("H6").Value = 0.139482
$ws.Range("I6").Value = 0.05587126560541624
$ws.Range("J6").Value = 0.05587126560541625
$ws.Range("M6").Value = 0.1629343333333333
$ws.Range("N6").Value = 0.488803
$ws.Range("O6").Value = 0.01139829922935564
$ws.Range("P6").Value = 0.01139829922935564
$ws.Range("Q6").Value = 0.007575468893999999
$ws.Range("R6").Value = 0.068179220046
$ws.Range("S6").Value = 0.0006368374036933405
$ws.Range("T6").Value = 0.0006368374036933406

# Row 7
$ws.Range("I7").Value = 0.7083039420562237
$ws.Range("J7").Value = 0.7083039420562237
$ws.Range("M7").Value = 0.7592313333333333
$ws.Range("N7").Value = 2.277694
$ws.Range("O7").Value = 0.05311309006881704
$ws.Range("P7").Value = 0.05311309006881704
$ws.Range("Q7").Value = 0.4475094224957777
$ws.Range("R7").Value = 4.027584802462
$ws.Range("S7").Value = 0.03762021107053037
$ws.Range("T7").Value = 0.03762021107053037

# Row 8
$ws.Range("I8").Value = 0.7083039420562237
$ws.Range("J8").Value = 0.7083039420562237
$ws.Range("O8").Value = 0.1721303313829956
$ws.Range("P8").Value = 0.1721303313829956
$ws.Range("S8").Value = 0.1219205922660199
$ws.Range("T8").Value = 0.1219205922660199

# Row 9
$ws.Range("I9").Value = 0.7083039420562237
$ws.Range("J9").Value = 0.7083039420562237
$ws.Range("M9").Value = 1.018760333333333
$ws.Range("N9").Value = 3.056281
$ws.Range("O9").Value = 0.07126880433834142
$ws.Range("P9").Value = 0.07126880433834142
$ws.Range("Q9").Value = 0.6004821303014444
$ws.Range("R9").Value = 5.404339172713
$ws.Range("S9").Value = 0.05047997505848092
$ws.Range("T9").Value = 0.05047997505848092

# Row 10
$ws.Range("I10").Value = 0.7083039420562237
$ws.Range("J10").Value = 0.7083039420562237
$ws.Range("M10").Value = 9.893154666666666
$ws.Range("N10").Value = 29.679464
$ws.Range("O10").Value = 0.6920894749804902
$ws.Range("P10").Value = 0.6920894749804903
$ws.Range("Q10").Value = 5.831266093963555
$ws.Range("R10").Value = 52.481394845672
$ws.Range("S10").Value = 0.4902097033843034
$ws.Range("T10").Value = 0.4902097033843035

# Row 11
$ws.Range("I11").Value = 0.7083039420562237
$ws.Range("J11").Value = 0.7083039420562237
$ws.Range("M11").Value = 0.1629343333333333
$ws.Range("N11").Value = 0.488803
$ws.Range("O11").Value = 0.01139829922935564
$ws.Range("P11").Value = 0.01139829922935564
$ws.Range("Q11").Value = 0.0960374608021111
$ws.Range("R11").Value = 0.864337147219
$ws.Range("S11").Value = 0.00807346027688902
$ws.Range("T11").Value = 0.00807346027688902

# Row 12
$ws.Range("G12").Value = 0.1246316666666667
$ws.Range("H12").Value = 0.373895
$ws.Range("I12").Value = 0.1497683346491813
$ws.Range("J12").Value = 0.1497683346491813
$ws.Range("M12").Value = 0.7592313333333333
$ws.Range("N12").Value = 2.277694
$ws.Range("O12").Value = 0.05311309006881704
$ws.Range("P12").Value = 0.05311309006881704
$ws.Range("Q12").Value = 0.09462426645888887
$ws.Range("R12").Value = 0.85161839813
$ws.Range("S12").Value = 0.007954659047678695
$ws.Range("T12").Value = 0.007954659047678697

# Row 13
$ws.Range("G13").Value = 0.1246316666666667
$ws.Range("H13").Value = 0.373895
$ws.Range("I13").Value = 0.1497683346491813
$ws.Range("J13").Value = 0.1497683346491813
$ws.Range("O13").Value = 0.1721303313829956
$ws.Range("P13").Value = 0.1721303313829956
$ws.Range("Q13").Value = 0.3066608687488889
$ws.Range("R13").Value = 2.75994781874
$ws.Range("S13").Value = 0.02577967307384295
$ws.Range("T13").Value = 0.02577967307384296

# Row 14
$ws.Range("G14").Value = 0.1246316666666667
$ws.Range("H14").Value = 0.373895
$ws.Range("I14").Value = 0.1497683346491813
$ws.Range("J14").Value = 0.1497683346491813
$ws.Range("M14").Value = 1.018760333333333
$ws.Range("N14").Value = 3.056281
$ws.Range("O14").Value = 0.07126880433834142
$ws.Range("P14").Value = 0.07126880433834142
$ws.Range("Q14").Value = 0.1269697982772222
$ws.Range("R14").Value = 1.142728184495
$ws.Range("S14").Value = 0.01067381013819174
$ws.Range("T14").Value = 0.01067381013819174

# Row 15
$ws.Range("G15").Value = 0.1246316666666667
$ws.Range("H15").Value = 0.373895
$ws.Range("I15").Value = 0.1497683346491813
$ws.Range("J15").Value = 0.1497683346491813
$ws.Range("M15").Value = 9.893154666666666
$ws.Range("N15").Value = 29.679464
$ws.Range("O15").Value = 0.6920894749804902
$ws.Range("P15").Value = 0.6920894749804903
$ws.Range("Q15").Value = 1.233000354697777
$ws.Range("R15").Value = 11.09700319228
$ws.Range("S15").Value = 0.1036530880960542
$ws.Range("T15").Value = 0.1036530880960543

# Row 16
$ws.Range("G16").Value = 0.1246316666666667
$ws.Range("H16").Value = 0.373895
$ws.Range("I16").Value = 0.1497683346491813
$ws.Range("J16").Value = 0.1497683346491813
$ws.Range("M16").Value = 0.1629343333333333
$ws.Range("N16").Value = 0.488803
$ws.Range("O16").Value = 0.01139829922935564
$ws.Range("P16").Value = 0.01139829922935564
$ws.Range("Q16").Value = 0.02030677752055555
$ws.Range("R16").Value = 0.182760997685
$ws.Range("S16").Value = 0.001707104293413641
$ws.Range("T16").Value = 0.001707104293413641

# Row 17
$ws.Range("G17").Value = 0.05954566666666666
$ws.Range("H17").Value = 0.178637
$ws.Range("I17").Value = 0.07155529225243931
$ws.Range("J17").Value = 0.07155529225243933
$ws.Range("M17").Value = 0.7592313333333333
$ws.Range("N17").Value = 2.277694
$ws.Range("O17").Value = 0.05311309006881704
$ws.Range("P17").Value = 0.05311309006881704
$ws.Range("Q17").Value = 0.04520893589755555
$ws.Range("R17").Value = 0.4068804230779999
$ws.Range("S17").Value = 0.003800522682304335
$ws.Range("T17").Value = 0.003800522682304336

# Row 18
$ws.Range("G18").Value = 0.05954566666666666
$ws.Range("H18").Value = 0.178637
$ws.Range("I18").Value = 0.07155529225243931
$ws.Range("J18").Value = 0.07155529225243933
$ws.Range("O18").Value = 0.1721303313829956
$ws.Range("P18").Value = 0.1721303313829956
$ws.Range("Q18").Value = 0.1465143358715555
$ws.Range("R18").Value = 1.318629022844
$ws.Range("S18").Value = 0.01231683616761948
$ws.Range("T18").Value = 0.01231683616761948

# Row 19
$ws.Range("G19").Value = 0.05954566666666666
$ws.Range("H19").Value = 0.178637
$ws.Range("I19").Value = 0.07155529225243931
$ws.Range("J19").Value = 0.07155529225243933
$ws.Range("M19").Value = 1.018760333333333
$ws.Range("N19").Value = 3.056281
$ws.Range("O19").Value = 0.07126880433834142
$ws.Range("P19").Value = 0.07126880433834142
$ws.Range("Q19").Value = 0.06066276322188888
$ws.Range("R19").Value = 0.5459648689969999
$ws.Range("S19").Value = 0.005099660122911935
$ws.Range("T19").Value = 0.005099660122911936

# Row 20
$ws.Range("G20").Value = 0.05954566666666666
$ws.Range("H20").Value = 0.178637
$ws.Range("I20").Value = 0.07155529225243931
$ws.Range("J20").Value = 0.07155529225243933
$ws.Range("M20").Value = 9.893154666666666
$ws.Range("N20").Value = 29.679464
$ws.Range("O20").Value = 0.6920894749804902
$ws.Range("P20").Value = 0.6920894749804903
$ws.Range("Q20").Value = 0.589094490063111
$ws.Range("R20").Value = 5.301850410568
$ws.Range("S20").Value = 0.04952266464706626
$ws.Range("T20").Value = 0.04952266464706628

# Row 21
$ws.Range("G21").Value = 0.05954566666666666
$ws.Range("H21").Value = 0.178637
$ws.Range("I21").Value = 0.07155529225243931
$ws.Range("J21").Value = 0.07155529225243933
$ws.Range("M21").Value = 0.1629343333333333
$ws.Range("N21").Value = 0.488803
$ws.Range("O21").Value = 0.01139829922935564
$ws.Range("P21").Value = 0.01139829922935564
$ws.Range("Q21").Value = 0.00970203350122222
$ws.Range("R21").Value = 0.087318301511
$ws.Range("S21").Value = 0.000815608632537297
$ws.Range("T21").Value = 0.0008156086325372971

# Row 22
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 0.3333333333333333
$ws.Range("G22").Value = 0.01206733333333333
$ws.Range("H22").Value = 0.036202
$ws.Range("I22").Value = 0.01450116543673935
$ws.Range("J22").Value = 0.01450116543673936
$ws.Range("M22").Value = 0.7592313333333333
$ws.Range("N22").Value = 2.277694
$ws.Range("O22").Value = 0.05311309006881704
$ws.Range("P22").Value = 0.05311309006881704
$ws.Range("Q22").Value = 0.009161897576444443
$ws.Range("R22").Value = 0.082457078188
$ws.Range("S22").Value = 0.0007702017059443538
$ws.Range("T22").Value = 0.000770201705944354

# Row 23
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 0.3333333333333333
$ws.Range("G23").Value = 0.01206733333333333
$ws.Range("H23").Value = 0.036202
$ws.Range("I23").Value = 0.01450116543673935
$ws.Range("J23").Value = 0.01450116543673936
$ws.Range("O23").Value = 0.1721303313829956
$ws.Range("P23").Value = 0.1721303313829956
$ws.Range("Q23").Value = 0.02969212418044444
$ws.Range("R23").Value = 0.267229117624
$ws.Range("S23").Value = 0.002496090412065587
$ws.Range("T23").Value = 0.002496090412065587

# Row 24
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 0.3333333333333333
$ws.Range("G24").Value = 0.01206733333333333
$ws.Range("H24").Value = 0.036202
$ws.Range("I24").Value = 0.01450116543673935
$ws.Range("J24").Value = 0.01450116543673936
$ws.Range("M24").Value = 1.018760333333333
$ws.Range("N24").Value = 3.056281
$ws.Range("O24").Value = 0.07126880433834142
$ws.Range("P24").Value = 0.07126880433834142
$ws.Range("Q24").Value = 0.01229372052911111
$ws.Range("R24").Value = 0.110643484762
$ws.Range("S24").Value = 0.001033480722188896
$ws.Range("T24").Value = 0.001033480722188896

# Row 25
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 0.3333333333333333
$ws.Range("G25").Value = 0.01206733333333333
$ws.Range("H25").Value = 0.036202
$ws.Range("I25").Value = 0.01450116543673935
$ws.Range("J25").Value = 0.01450116543673936
$ws.Range("M25").Value = 9.893154666666666
$ws.Range("N25").Value = 29.679464
$ws.Range("O25").Value = 0.6920894749804902
$ws.Range("P25").Value = 0.6920894749804903
$ws.Range("Q25").Value = 0.1193839950808889
$ws.Range("R25").Value = 1.074455955728
$ws.Range("S25").Value = 0.01003610397371817
$ws.Range("T25").Value = 0.01003610397371817

# Row 26
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 0.3333333333333333
$ws.Range("G26").Value = 0.01206733333333333
$ws.Range("H26").Value = 0.036202
$ws.Range("I26").Value = 0.01450116543673935
$ws.Range("J26").Value = 0.01450116543673936
$ws.Range("M26").Value = 0.1629343333333333
$ws.Range("N26").Value = 0.488803
$ws.Range("O26").Value = 0.01139829922935564
$ws.Range("P26").Value = 0.01139829922935564
$ws.Range("Q26").Value = 0.001966182911777777
$ws.Range("R26").Value = 0.017695646206
$ws.Range("S26").Value = 0.0001652886228223449
$ws.Range("T26").Value = 0.0001652886228223449
